# Insert three new runs ("$\SunQuarTeX$", " ", "Example - ") at the
# very start of the first paragraph (style "a5"), before the existing
# "课表" run. The new runs carry no run-properties, so we build them as
# a raw OOXML fragment and insert it via Range.InsertXML - this avoids
# Word's normal "inherit the neighbouring run's formatting" behaviour
# that InsertBefore/InsertAfter would otherwise apply (e.g. picking up
# the <w:rFonts w:hint="eastAsia"/> from the "课表" run).

$d = $word.ActiveDocument

$target = $d.Paragraphs(1).Range
$insertionPoint = $target.Duplicate
$insertionPoint.Collapse(1)   ## wdCollapseStart

$openXmlFragment = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r><w:t xml:space="preserve">$\SunQuarTeX$</w:t></w:r>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:r><w:t xml:space="preserve">Example - </w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insertionPoint.InsertXML($openXmlFragment)
